# Auto-generated edit script applying the Titan_Profits profit-recalculation diff
# to the underlying ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 683517
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 717672.9
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 2153018.7
$ws.Range("M17").Value = -1032
$ws.Range("N17").Value = -2153354.7
# Row 19
$ws.Range("H19").Value = 1365.6428
$ws.Range("I19").Value = 997.5
$ws.Range("K19").Value = 997.5
$ws.Range("M19").Value = -822.5
# Row 106
$ws.Range("H106").Value = 3369786
$ws.Range("I106").Value = 3706411.2
$ws.Range("K106").Value = 3706411.2
$ws.Range("M106").Value = -3705780.2
# Row 112
$ws.Range("H112").Value = 8022287.5
$ws.Range("J112").Value = 10490503
$ws.Range("L112").Value = 31471509
$ws.Range("N112").Value = -31473725
# Row 132
$ws.Range("H132").Value = 249792.1
$ws.Range("I132").Value = 270746.66
$ws.Range("J132").Value = 61201.2
$ws.Range("K132").Value = 812239.98
$ws.Range("L132").Value = 183603.6
$ws.Range("M132").Value = -809709.98
$ws.Range("N132").Value = -188663.6
# Row 138
$ws.Range("H138").Value = 2621.5684
$ws.Range("I138").Value = 1902.5454
$ws.Range("J138").Value = 3241.9019
$ws.Range("K138").Value = 5707.6362
$ws.Range("L138").Value = 9725.705699999999
$ws.Range("M138").Value = -567.6361999999999
$ws.Range("N138").Value = -20005.7057

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17649.045
$ws.Range("I32").Value = 2524.1587
$ws.Range("K32").Value = 2524.1587
$ws.Range("M32").Value = -2237.1587
# Row 61
$ws.Range("H61").Value = 2034.3478
$ws.Range("I61").Value = 1399.6666
$ws.Range("K61").Value = 1399.6666
$ws.Range("M61").Value = -1187.6666
# Row 74
$ws.Range("H74").Value = 3629.7708
$ws.Range("I74").Value = 900.36584
$ws.Range("J74").Value = 19616.285
$ws.Range("K74").Value = 900.36584
$ws.Range("L74").Value = 19616.285
$ws.Range("M74").Value = -26.36584000000005
$ws.Range("N74").Value = -21364.285
# Row 77
$ws.Range("H77").Value = 3629.7708
$ws.Range("I77").Value = 900.36584
$ws.Range("J77").Value = 19616.285
$ws.Range("K77").Value = 4501.8292
$ws.Range("L77").Value = 98081.425
$ws.Range("M77").Value = -133.8292000000001
$ws.Range("N77").Value = -106817.425
# Row 97
$ws.Range("H97").Value = 37048736
$ws.Range("I97").Value = 37048736
$ws.Range("K97").Value = 37048736
$ws.Range("M97").Value = -37048240
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
# Row 132
$ws.Range("H132").Value = 2499.2917
$ws.Range("I132").Value = 2096.9092
$ws.Range("K132").Value = 6290.7276
$ws.Range("M132").Value = -3760.7276
# Row 136
$ws.Range("H136").Value = 2034.3478
$ws.Range("I136").Value = 1399.6666
$ws.Range("K136").Value = 4198.9998
$ws.Range("M136").Value = -1648.9998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -227
# Row 94
$ws.Range("H94").Value = 1096.5
$ws.Range("I94").Value = 858.4706
$ws.Range("J94").Value = 1905.8
$ws.Range("K94").Value = 858.4706
$ws.Range("L94").Value = 1905.8
$ws.Range("M94").Value = -407.4706
$ws.Range("N94").Value = -2807.8
# Row 132
$ws.Range("H132").Value = 36114.668
$ws.Range("J132").Value = 36114.668
$ws.Range("L132").Value = 36114.668
$ws.Range("N132").Value = -46234.668
# Row 134
$ws.Range("H134").Value = 2261.1562
$ws.Range("I134").Value = 1443.8723
$ws.Range("J134").Value = 4520.706
$ws.Range("K134").Value = 4331.6169
$ws.Range("L134").Value = 13562.118
$ws.Range("M134").Value = -1796.6169
$ws.Range("N134").Value = -18632.118

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2020
$ws.Range("I31").Value = 1085.3334
$ws.Range("K31").Value = 1085.3334
$ws.Range("M31").Value = -790.3334
# Row 34
$ws.Range("H34").Value = 2020
$ws.Range("I34").Value = 1085.3334
$ws.Range("K34").Value = 1085.3334
$ws.Range("M34").Value = -883.3334
# Row 134
$ws.Range("H134").Value = 2335.1277
$ws.Range("I134").Value = 1563.4849
$ws.Range("J134").Value = 4154
$ws.Range("K134").Value = 4690.4547
$ws.Range("L134").Value = 12462
$ws.Range("M134").Value = -2155.4547
$ws.Range("N134").Value = -17532

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 616.73914
$ws.Range("I122").Value = 314.53845
$ws.Range("J122").Value = 1009.6
$ws.Range("K122").Value = 2830.84605
$ws.Range("L122").Value = 9086.4
$ws.Range("M122").Value = -380.8460500000001
$ws.Range("N122").Value = -13986.4
# Row 132
$ws.Range("H132").Value = 13889643
$ws.Range("I132").Value = 772
$ws.Range("J132").Value = 33334062
$ws.Range("K132").Value = 6948
$ws.Range("L132").Value = 300006558
$ws.Range("M132").Value = -4418
$ws.Range("N132").Value = -300011618

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 1666.6666
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -888
$ws.Range("N7").Value = -3224
# Row 8
$ws.Range("H8").Value = 1666.6666
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = -861
$ws.Range("N8").Value = -3278
# Row 97
$ws.Range("H97").Value = 46182.547
$ws.Range("I97").Value = 53195.05
$ws.Range("K97").Value = 53195.05
$ws.Range("M97").Value = -52699.05
# Row 126
$ws.Range("H126").Value = 2254.3
$ws.Range("I126").Value = 1958.4
$ws.Range("J126").Value = 2352.9333
$ws.Range("K126").Value = 5875.200000000001
$ws.Range("L126").Value = 7058.7999
$ws.Range("M126").Value = -3405.200000000001
$ws.Range("N126").Value = -11998.7999
# Row 132
$ws.Range("H132").Value = 2206.4443
$ws.Range("I132").Value = 2100.56
$ws.Range("J132").Value = 2447.0908
$ws.Range("K132").Value = 6301.68
$ws.Range("L132").Value = 7341.2724
$ws.Range("M132").Value = -3771.68
$ws.Range("N132").Value = -12401.2724
# Row 135
$ws.Range("H135").Value = 142891710
$ws.Range("J135").Value = 142891710
$ws.Range("L135").Value = 142891710
$ws.Range("N135").Value = -142901850

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 14075
$ws.Range("I14").Value = 17766.666
$ws.Range("J14").Value = 3000
$ws.Range("K14").Value = 17766.666
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = -17594.666
$ws.Range("N14").Value = -3344
# Row 16
$ws.Range("H16").Value = 1115.2106
$ws.Range("I16").Value = 1246.3334
$ws.Range("J16").Value = 623.5
$ws.Range("K16").Value = 1246.3334
$ws.Range("L16").Value = 623.5
$ws.Range("M16").Value = -1076.3334
# Row 40
$ws.Range("H40").Value = 5750
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 5750
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 5750
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -6022
# Row 122
$ws.Range("H122").Value = 3643.2144
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3643.2144
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 10929.6432
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -15829.6432
# Row 136
$ws.Range("H136").Value = 3476.9434
$ws.Range("I136").Value = 2068.7673
$ws.Range("J136").Value = 9532.1
$ws.Range("K136").Value = 6206.3019
$ws.Range("L136").Value = 28596.3
$ws.Range("M136").Value = -3656.3019
$ws.Range("N136").Value = -33696.3

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
# Row 126
$ws.Range("H126").Value = 92318.27
$ws.Range("I126").Value = 101350.1
$ws.Range("K126").Value = 304050.3
$ws.Range("M126").Value = -301580.3
# Row 136
$ws.Range("H136").Value = 7961531
$ws.Range("I136").Value = 8799235
$ws.Range("J136").Value = 3341.25
$ws.Range("K136").Value = 26397705
$ws.Range("L136").Value = 10023.75
$ws.Range("M136").Value = -26395155
$ws.Range("N136").Value = -15123.75
